$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.287.40"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "3.524.25"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.06"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.23"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("E9").Value = "  +7.47%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "4.136.98"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.66"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "67.201.08"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "3.526.98"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.35"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.15"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.82"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.98"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.45"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000123"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.24"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.181"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.29"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.07"
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  +5.01%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.897"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.91"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.72"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0746"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.52"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.22"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("E43").Value = "  +4.56%  "
$ws.Range("D44").Value = "2.811.53"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.96"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.94"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.52"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.852"
$ws.Range("E51").Value = "  +0.41%  "
